# Apply the "prior-authorization-indicator" StructureDefinition spreadsheet
# refresh: IBM/Alvearie branding -> LinuxForHealth branding, version bump,
# regenerated date, and the FHIR IG Publisher's corrected placement of the
# ele-1/ext-1 constraint text (it belongs on Extension.extension's row, not
# the root Extension row).

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet (Property / Value table) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/prior-authorization-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet (element definitions table) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" (root element): the ele-1/ext-1 constraint text was
# mis-placed here; the regenerated IG clears it from the root row.
$elements.Range("AI2").Value = ""

# Row 5 = "Extension.url": its Fixed Value mirrors the structure's own URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/prior-authorization-indicator"
